# Consolidate text runs when possible (mirrors the PowerPoint-writer change
# that merges adjacent <a:r> runs holding a single word/space each into one
# run per paragraph).
#
# The COM text-setter in this runtime keeps the *existing* run(s) (and their
# <a:rPr/>) when the new string is no longer than the current text, only
# dropping/merging the runs that become superfluous; it only fabricates a
# brand new run (with a "lang" attribute) when the text grows past what the
# current runs can hold. So, to end up with a single pre-existing run (and
# therefore no extra "lang" attribute) we first shrink the text down to a
# short prefix - which collapses everything onto the first run - and then
# grow it back out to the real target text.

$p = $ppt.ActivePresentation

function Consolidate-Text {
    param($shape, [string]$finalText)

    $tr = $shape.TextFrame.TextRange
    # Shrink to a 1-character prefix of the final text first so every run
    # but the first is dropped and the first run's <a:rPr/> is preserved.
    $tr.Text = $finalText.Substring(0, 1)
    $tr.Text = $finalText
}

# --- Slide titles: "Slide" " " "N" -> "Slide N" -------------------------

Consolidate-Text $p.Slides.Item(2).Shapes.Title "Slide 1"
Consolidate-Text $p.Slides.Item(4).Shapes.Title "Slide 3"
Consolidate-Text $p.Slides.Item(5).Shapes.Title "Slide 4"
Consolidate-Text $p.Slides.Item(6).Shapes.Title "Slide 5"

# --- Notes page for slide 7: word-by-word runs -> single run ------------

$notesShape = $p.Slides.Item(7).NotesPage.Shapes.Item(2)
$notesShape.TextFrame.TextRange.Text = "This is a blank slide: does it have a footer?"
